# netCrypto.xlsx edit:
#  - Row 3 (Deposit/Crypto/ETH/980.7002) was removed, shifting every row
#    below it up by one.
#  - Three new transaction rows were appended at the bottom of the table.
#  - The active selection/scroll position was changed to E3:T9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the old row 3 - everything below shifts up by one row.
$ws.Rows(3).Delete()

# 2) Append the new transaction rows at the bottom (now rows 159-161).
$ws.Range("E159").Value2 = "Deposit"
$ws.Range("N159").Value2 = "Crypto"
$ws.Range("P159").Value2 = "ETH"
$ws.Range("T159").Value2 = 2173.6392000000001

$ws.Range("E160").Value2 = "Deposit"
$ws.Range("N160").Value2 = "Credit Card"
$ws.Range("P160").Value2 = "Startrading"
$ws.Range("T160").Value2 = 267.1275

$ws.Range("E161").Value2 = "Withdrawal"
$ws.Range("N161").Value2 = "Credit Card"
$ws.Range("P161").Value2 = "Tradeprof"
$ws.Range("T161").Value2 = 267.13

# 3) Update the view: scroll/selection now highlights E3:T9.
$ws.Activate()
$ws.Range("E3:T9").Select()
